$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix "Time in session" values: correct spacing/punctuation in the
# hours/minutes strings (B7, C7).
$ws.Range("B7").Value = "  1,278 hrs., 15'"
$ws.Range("C7").Value = "829 hrs., 11'"

# Fix the "Measures passed, total" House value: replace the placeholder
# text "n2" with the real numeric value 772.
$ws.Range("C13").Value = 772
